$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Part 1: "Within rust belt regions, " / "if the regions are swing
# states" / ", Democrats link ... (worker, " used to be three separate
# runs with the _GoBack bookmark sitting between the first two. Those
# three runs collapse into a single run (the text itself is unchanged)
# and the _GoBack bookmark moves away from here entirely.
# ---------------------------------------------------------------------

$bm1 = $d.Bookmarks("_GoBack")
$bm1.Delete()

$rng1 = $d.Content
$rng1.Find.Execute("Within rust belt regions, if the regions are swing states, Democrats link threatened job security with domestic words (worker, ")
$start1 = $rng1.Start
$end1 = $rng1.End

# Pin the boundary right before "Within" with a temporary bookmark so
# that the upcoming text assignment can't coalesce backwards into the
# preceding "H3-2: " run.
$d.Bookmarks.Add("ZZPIN1", $d.Range($start1, $start1))

# Force a genuine text mutation (so the engine actually rebuilds the
# run) by writing the text with a harmless marker character inserted,
# then stripping that marker back out again. A pure no-op assignment
# would leave the original run boundaries untouched.
$markedText1 = "Within rust belt regions, if the regions are swing states~, Democrats link threatened job security with domestic words (worker, "
$d.Range($start1, $end1).Text = $markedText1

$markerPos1 = $start1 + "Within rust belt regions, if the regions are swing states".Length
$d.Range($markerPos1, $markerPos1 + 1).Text = ""

$d.Bookmarks("ZZPIN1").Delete()

# ---------------------------------------------------------------------
# Part 2: "author's position as the incumbent vs. challenger" becomes
# "candidate from incumbent party vs. opposing party", and the
# _GoBack bookmark reappears right after it (before "the partisan
# similarity ...").
# ---------------------------------------------------------------------

$rng2 = $d.Content
$rng2.Find.Execute("author’s position as the incumbent vs. challenger, ")
$start2 = $rng2.Start
$end2 = $rng2.End

$d.Range($start2, $end2).Text = ""

$cursor2 = $d.Range($start2, $start2)
$cursor2.InsertAfter("candidate from incumbent party vs.")
$cursor2.Collapse(0)
$cursor2.InsertAfter(" opposing party")
$cursor2.Collapse(0)
$cursor2.InsertAfter(",")
$cursor2.Collapse(0)
$cursor2.InsertAfter(" ")
$cursor2.Collapse(0)

$d.Bookmarks.Add("_GoBack", $cursor2)
